# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (previously the "latest" row) loses its special date-only format
# and becomes a normal date+time formatted row.
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 10 becomes the latest entry, with the special date-only format.
$ws.Range("A10").Value = 45959
$ws.Range("A10").NumberFormat = "YYYY-MM-DD"
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = 22
